$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.004488110542297
$ws.Range("B1").Value = 1.494637370109558
$ws.Range("C1").Value = -1
$ws.Range("D1").Value = -1
$ws.Range("E1").Value = 0.7930213212966919
